$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E1:E30").Copy() | Out-Null
$ws.Range("F1:F30").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = 42155
$ws.Range("F3").Value = 6010000
$ws.Range("F4").Value = 412000
$ws.Range("F5").Value = 3011000
$ws.Range("F6").Value = 2587000
$ws.Range("F7").Value = 15587000
$ws.Range("F8").Value = 4337000
$ws.Range("F9").Value = 3358000
$ws.Range("F10").Value = 7892000
$ws.Range("F11").Value = 5924000
$ws.Range("F12").Value = 21597000
$ws.Range("F14").Value = 12707000
$ws.Range("F15").Value = 3000
$ws.Range("F16").Value = 12704000
$ws.Range("F17").Value = 2558000
$ws.Range("F18").Value = 1079000
$ws.Range("F19").Value = 1479000
$ws.Range("F20").Value = "n.a."
$ws.Range("F21").Value = 6332000
$ws.Range("F22").Value = 107000
$ws.Range("F23").Value = 2131000
$ws.Range("F24").Value = 4094000
$ws.Range("F25").Value = 21597000
$ws.Range("F27").Value = 5564000
$ws.Range("F28").Value = 9255000
$ws.Range("F29").Value = 64694223.414803401
$ws.Range("F30").Value = 62600

$ws.Range("A1:F30").ColumnWidth = 14.91
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.Range("J11").Select() | Out-Null
Write-Host "done"
